$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8972259163856506
$ws.Range("B1").Value = 1.306082129478455
$ws.Range("C1").Value = 5.2569580078125
$ws.Range("D1").Value = 1.611467838287354
$ws.Range("E1").Value = 0.9421895742416382
